# Add a team Win/Loss/Tie record alongside the existing player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constant values (COM enums aren't exposed as bare names here):
#   xlCenter = -4108, xlTop = -4160, xlContinuous = 1, xlThin = 2
#   xlEdgeLeft = 7, xlEdgeTop = 8, xlEdgeBottom = 9, xlEdgeRight = 10
$xlCenter = -4108
$xlTop = -4160
$xlContinuous = 1
$xlThin = 2
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

# New header cells: AD1 = Wins, AE1 = Losses, AF1 = Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting already used by the other header cells (A1:AC1):
# bold, centered horizontally, aligned to top, thin border all the way around.
foreach ($addr in @("AD1", "AE1", "AF1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlTop

    $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeLeft).Weight = $xlThin
    $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeTop).Weight = $xlThin
    $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeBottom).Weight = $xlThin
    $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeRight).Weight = $xlThin
}

# Every player row (2-49) gets the same team record: 76 wins, 86 losses, 0 ties.
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 76
    $ws.Cells.Item($r, 31).Value = 86
    $ws.Cells.Item($r, 32).Value = 0
}
